$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Sheet2" positioned after the existing Sheet1
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate the summary details about process completion
$ws2.Range("A1").Value = "Total Process"
$ws2.Range("B1").Value = 100

$ws2.Range("A2").Value = "AWT (Average Waiting Time)"
$ws2.Range("B2").Value = 1172.99

$ws2.Range("A3").Value = "Total Waiting Time"
$ws2.Range("B3").Value = 117299

$ws2.Range("A4").Value = "ATAT (Average Turn Around Time)"
$ws2.Range("B4").Value = 1185.57

$ws2.Range("A5").Value = "Total Turn Around Time"
$ws2.Range("B5").Value = 118557
